# Auto-generated edit script applying scheduled market-data refresh
# to the Gilgamesh_Profits leve-profit workbook (per-sheet H:N recompute).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2028.8572
$ws.Range("I43").Value = 2566.6667
$ws.Range("K43").Value = 2566.6667
$ws.Range("M43").Value = -2497.6667
$ws.Range("H62").Value = 3435.5454
$ws.Range("I62").Value = 3285.375
$ws.Range("J62").Value = 3836
$ws.Range("K62").Value = 3285.375
$ws.Range("L62").Value = 3836
$ws.Range("M62").Value = -2661.375
$ws.Range("N62").Value = -5084
$ws.Range("H65").Value = 3435.5454
$ws.Range("I65").Value = 3285.375
$ws.Range("J65").Value = 3836
$ws.Range("K65").Value = 16426.875
$ws.Range("L65").Value = 19180
$ws.Range("M65").Value = -13306.875
$ws.Range("N65").Value = -25420
$ws.Range("H112").Value = 2060.2693
$ws.Range("J112").Value = 2102.7917
$ws.Range("L112").Value = 6308.375100000001
$ws.Range("N112").Value = -8524.375100000001
$ws.Range("H132").Value = 8736.936
$ws.Range("I132").Value = 6959.3335
$ws.Range("K132").Value = 20878.0005
$ws.Range("M132").Value = -18348.0005
$ws.Range("H137").Value = 1520680.9
$ws.Range("I137").Value = 2632609.8
$ws.Range("K137").Value = 7897829.399999999
$ws.Range("M137").Value = -7895279.399999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 591.8
$ws.Range("I2").Value = 299.5
$ws.Range("K2").Value = 299.5
$ws.Range("M2").Value = -186.5
$ws.Range("I5").Value = 249
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 249
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -137
$ws.Range("N5").ClearContents()
$ws.Range("H32").Value = 4716.095
$ws.Range("I32").Value = 4205.9707
$ws.Range("K32").Value = 4205.9707
$ws.Range("M32").Value = -3918.9707
$ws.Range("H45").Value = 32213.143
$ws.Range("I45").Value = 36579.75
$ws.Range("K45").Value = 36579.75
$ws.Range("M45").Value = -36202.75
$ws.Range("H61").Value = 2468.0967
$ws.Range("I61").Value = 1721.75
$ws.Range("J61").Value = 3264.2
$ws.Range("K61").Value = 1721.75
$ws.Range("L61").Value = 3264.2
$ws.Range("M61").Value = -1509.75
$ws.Range("N61").Value = -3688.2
$ws.Range("H63").Value = 1783.2222
$ws.Range("J63").Value = 1951
$ws.Range("L63").Value = 1951
$ws.Range("N63").Value = -3323
$ws.Range("H66").Value = 1783.2222
$ws.Range("J66").Value = 1951
$ws.Range("L66").Value = 9755
$ws.Range("N66").Value = -16619
$ws.Range("H74").Value = 165816.44
$ws.Range("I74").Value = 327957.6
$ws.Range("K74").Value = 327957.6
$ws.Range("M74").Value = -327083.6
$ws.Range("H77").Value = 165816.44
$ws.Range("I77").Value = 327957.6
$ws.Range("K77").Value = 1639788
$ws.Range("M77").Value = -1635420
$ws.Range("H116").Value = 591.8
$ws.Range("I116").Value = 299.5
$ws.Range("K116").Value = 299.5
$ws.Range("M116").Value = 1994.5
$ws.Range("H132").Value = 2000.8158
$ws.Range("I132").Value = 1401.0333
$ws.Range("K132").Value = 4203.0999
$ws.Range("M132").Value = -1673.0999
$ws.Range("H136").Value = 2468.0967
$ws.Range("I136").Value = 1721.75
$ws.Range("J136").Value = 3264.2
$ws.Range("K136").Value = 5165.25
$ws.Range("L136").Value = 9792.599999999999
$ws.Range("M136").Value = -2615.25
$ws.Range("N136").Value = -14892.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 591.8
$ws.Range("I3").Value = 299.5
$ws.Range("K3").Value = 299.5
$ws.Range("M3").Value = -185.5
$ws.Range("I4").Value = 249
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 249
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -134
$ws.Range("N4").ClearContents()
$ws.Range("H110").Value = 50000
$ws.Range("J110").Value = 50000
$ws.Range("L110").Value = 50000
$ws.Range("N110").Value = -58180
$ws.Range("H134").Value = 3899.1794
$ws.Range("I134").Value = 4150
$ws.Range("K134").Value = 12450
$ws.Range("M134").Value = -9915
$ws.Range("H135").Value = 98983.47
$ws.Range("J135").Value = 98983.47
$ws.Range("L135").Value = 98983.47
$ws.Range("N135").Value = -109123.47
$ws.Range("H138").Value = 64797.273
$ws.Range("J138").Value = 64797.273
$ws.Range("L138").Value = 64797.273
$ws.Range("N138").Value = -75077.273

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 200000060
$ws.Range("J7").Value = 500000030
$ws.Range("L7").Value = 500000030
$ws.Range("N7").Value = -500000256
$ws.Range("H22").Value = 2396
$ws.Range("I22").Value = 2396
$ws.Range("K22").Value = 2396
$ws.Range("M22").Value = -2046
$ws.Range("H31").Value = 4504.615
$ws.Range("I31").Value = 3425
$ws.Range("K31").Value = 3425
$ws.Range("M31").Value = -3130
$ws.Range("H34").Value = 4504.615
$ws.Range("I34").Value = 3425
$ws.Range("K34").Value = 3425
$ws.Range("M34").Value = -3223
$ws.Range("H58").Value = 2615.8076
$ws.Range("I58").Value = 1670.0834
$ws.Range("K58").Value = 1670.0834
$ws.Range("M58").Value = -1467.0834
$ws.Range("H107").Value = 606.7143
$ws.Range("I107").Value = 455.33334
$ws.Range("K107").Value = 455.33334
$ws.Range("M107").Value = 1464.66666
$ws.Range("H132").Value = 2499.3076
$ws.Range("I132").Value = 1812.125
$ws.Range("K132").Value = 5436.375
$ws.Range("M132").Value = -2906.375
$ws.Range("H133").Value = 71996.664
$ws.Range("J133").Value = 71996.664
$ws.Range("L133").Value = 71996.664
$ws.Range("N133").Value = -77056.664
$ws.Range("H136").Value = 2615.8076
$ws.Range("I136").Value = 1670.0834
$ws.Range("K136").Value = 5010.2502
$ws.Range("M136").Value = -2460.2502

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1350.4375
$ws.Range("I5").Value = 350.2857
$ws.Range("J5").Value = 2128.3333
$ws.Range("K5").Value = 1050.8571
$ws.Range("L5").Value = 6384.999899999999
$ws.Range("M5").Value = -938.8571000000002
$ws.Range("N5").Value = -6608.999899999999
$ws.Range("H50").Value = 611
$ws.Range("I50").Value = 750
$ws.Range("J50").Value = 333
$ws.Range("K50").Value = 2250
$ws.Range("L50").Value = 999
$ws.Range("N50").Value = -1961
$ws.Range("M50").Value = -1769
$ws.Range("H53").Value = 611
$ws.Range("I53").Value = 750
$ws.Range("J53").Value = 333
$ws.Range("K53").Value = 2250
$ws.Range("L53").Value = 999
$ws.Range("N53").Value = -1961
$ws.Range("M53").Value = -1769
$ws.Range("H135").Value = 1350.4375
$ws.Range("I135").Value = 350.2857
$ws.Range("J135").Value = 2128.3333
$ws.Range("K135").Value = 3152.5713
$ws.Range("L135").Value = 19154.9997
$ws.Range("M135").Value = -617.5713000000001
$ws.Range("N135").Value = -24224.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 3967
$ws.Range("I97").Value = 3348
$ws.Range("J97").Value = 4338.4
$ws.Range("K97").Value = 3348
$ws.Range("L97").Value = 4338.4
$ws.Range("M97").Value = -2852
$ws.Range("N97").Value = -5330.4
$ws.Range("H113").Value = 4161.2085
$ws.Range("I113").Value = 3527
$ws.Range("J113").Value = 5218.222
$ws.Range("K113").Value = 3527
$ws.Range("L113").Value = 5218.222
$ws.Range("M113").Value = -1357
$ws.Range("N113").Value = -9558.222
$ws.Range("H132").Value = 2205.8918
$ws.Range("I132").Value = 1918.4482
$ws.Range("J132").Value = 3247.875
$ws.Range("K132").Value = 5755.3446
$ws.Range("L132").Value = 9743.625
$ws.Range("M132").Value = -3225.3446
$ws.Range("N132").Value = -14803.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H22").Value = 634.9
$ws.Range("I22").Value = 550
$ws.Range("K22").Value = 550
$ws.Range("M22").Value = -255
$ws.Range("H27").Value = 634.9
$ws.Range("I27").Value = 550
$ws.Range("K27").Value = 550
$ws.Range("M27").Value = -443
$ws.Range("H46").Value = 3000
$ws.Range("I46").Value = 3250
$ws.Range("K46").Value = 3250
$ws.Range("M46").Value = -3062
$ws.Range("H61").Value = 57002.5
$ws.Range("I61").Value = 4000
$ws.Range("J61").Value = 110005
$ws.Range("K61").Value = 4000
$ws.Range("L61").Value = 110005
$ws.Range("M61").Value = -3798
$ws.Range("N61").Value = -110409
$ws.Range("H113").Value = 57002.5
$ws.Range("I113").Value = 4000
$ws.Range("J113").Value = 110005
$ws.Range("K113").Value = 4000
$ws.Range("L113").Value = 110005
$ws.Range("M113").Value = -1830
$ws.Range("N113").Value = -114345
$ws.Range("H134").Value = 103990.664
$ws.Range("J134").Value = 103990.664
$ws.Range("L134").Value = 103990.664
$ws.Range("N134").Value = -114130.664

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 5000
$ws.Range("I8").Value = 5000
$ws.Range("J8").Value = 5000
$ws.Range("K8").Value = 5000
$ws.Range("L8").Value = 5000
$ws.Range("M8").Value = -4860
$ws.Range("N8").Value = -5280
$ws.Range("H113").Value = 490.41666
$ws.Range("I113").Value = 297
$ws.Range("K113").Value = 891
$ws.Range("M113").Value = 1279
$ws.Range("H140").Value = 114845.75
$ws.Range("J140").Value = 109664.336
$ws.Range("L140").Value = 109664.336
$ws.Range("N140").Value = -120024.336
